{"js": "// Update the title date line and the 20x5 arithmetic-problem table.\n// The title paragraph's single run text is replaced, and the table's\n// cell values are replaced in place (row-major order), which preserves\n// each cell's existing run formatting (font/size) automatically.\n\nconst body = context.document.body;\n\n// 1) Title paragraph: \"2025-12-11 Thursday\" -> \"2025-12-12 Friday\"\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text === \"2025-12-11 Thursday\") {\n  titlePara.getRange().insertText(\"2025-12-12 Friday\", Word.InsertLocation.replace);\n}\n\n// 2) Table of arithmetic problems: replace every cell's text (row-major).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"93-56=\", \"7+36=\", \"98-79=\", \"54+27=\", \"64+27=\"],\n  [\"58+15=\", \"90-51=\", \"7+47=\", \"9+8=\", \"92-73=\"],\n  [\"13+29=\", \"9+43=\", \"69+16=\", \"36+26=\", \"8+44=\"],\n  [\"43-26=\", \"14+17=\", \"66+29=\", \"54+17=\", \"93-49=\"],\n  [\"68+5=\", \"16+66=\", \"44+19=\", \"37+8=\", \"24+58=\"],\n  [\"20-1=\", \"56+7=\", \"90-42=\", \"39+33=\", \"67-49=\"],\n  [\"28+9=\", \"76+6=\", \"82-18=\", \"33+18=\", \"44+29=\"],\n  [\"95-68=\", \"37+39=\", \"93-55=\", \"84-48=\", \"97-8=\"],\n  [\"74-7=\", \"7+68=\", \"64+27=\", \"47+26=\", \"14-8=\"],\n  [\"13+19=\", \"29+55=\", \"8+24=\", \"43-8=\", \"78+15=\"],\n  [\"15-6=\", \"37+8=\", \"9+84=\", \"7+89=\", \"4+69=\"],\n  [\"73-38=\", \"19+39=\", \"56-27=\", \"72-54=\", \"84-48=\"],\n  [\"49+14=\", \"31-24=\", \"84-17=\", \"50-6=\", \"50-16=\"],\n  [\"27+35=\", \"78+14=\", \"35-29=\", \"30-27=\", \"91-73=\"],\n  [\"91-45=\", \"41-33=\", \"97-18=\", \"83-69=\", \"18+15=\"],\n  [\"19+68=\", \"29+42=\", \"73-8=\", \"44+38=\", \"37+25=\"],\n  [\"15+27=\", \"96-19=\", \"70-6=\", \"64+19=\", \"46+17=\"],\n  [\"43+48=\", \"61-5=\", \"79+5=\", \"38+7=\", \"43-16=\"],\n  [\"77+4=\", \"42+19=\", \"93-7=\", \"62-35=\", \"88+3=\"],\n  [\"75+17=\", \"45-16=\", \"83-37=\", \"17+68=\", \"16+16=\"]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the title date line and the 20x5 arithmetic-problem table.\n# The title paragraph's range text is replaced directly, and every table\n# cell's range text is replaced in place (row-major order), which keeps\n# each cell's existing run formatting (font/size) untouched.\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph: \"2025-12-11 Thursday\" -> \"2025-12-12 Friday\"\n# (Paragraph.Range.Text includes the trailing paragraph-mark char, so trim\n# it before comparing; assigning back to .Text replaces just the content.)\n$titleRange = $d.Paragraphs.Item(1).Range\n$titleText = $titleRange.Text.TrimEnd([char]13, [char]7)\nif ($titleText -eq \"2025-12-11 Thursday\") {\n    $titleRange.Text = \"2025-12-12 Friday\"\n}\n\n# 2) Table of arithmetic problems: replace every cell's text (row-major).\n$newValues = @(\n    @(\"93-56=\", \"7+36=\", \"98-79=\", \"54+27=\", \"64+27=\"),\n    @(\"58+15=\", \"90-51=\", \"7+47=\", \"9+8=\", \"92-73=\"),\n    @(\"13+29=\", \"9+43=\", \"69+16=\", \"36+26=\", \"8+44=\"),\n    @(\"43-26=\", \"14+17=\", \"66+29=\", \"54+17=\", \"93-49=\"),\n    @(\"68+5=\", \"16+66=\", \"44+19=\", \"37+8=\", \"24+58=\"),\n    @(\"20-1=\", \"56+7=\", \"90-42=\", \"39+33=\", \"67-49=\"),\n    @(\"28+9=\", \"76+6=\", \"82-18=\", \"33+18=\", \"44+29=\"),\n    @(\"95-68=\", \"37+39=\", \"93-55=\", \"84-48=\", \"97-8=\"),\n    @(\"74-7=\", \"7+68=\", \"64+27=\", \"47+26=\", \"14-8=\"),\n    @(\"13+19=\", \"29+55=\", \"8+24=\", \"43-8=\", \"78+15=\"),\n    @(\"15-6=\", \"37+8=\", \"9+84=\", \"7+89=\", \"4+69=\"),\n    @(\"73-38=\", \"19+39=\", \"56-27=\", \"72-54=\", \"84-48=\"),\n    @(\"49+14=\", \"31-24=\", \"84-17=\", \"50-6=\", \"50-16=\"),\n    @(\"27+35=\", \"78+14=\", \"35-29=\", \"30-27=\", \"91-73=\"),\n    @(\"91-45=\", \"41-33=\", \"97-18=\", \"83-69=\", \"18+15=\"),\n    @(\"19+68=\", \"29+42=\", \"73-8=\", \"44+38=\", \"37+25=\"),\n    @(\"15+27=\", \"96-19=\", \"70-6=\", \"64+19=\", \"46+17=\"),\n    @(\"43+48=\", \"61-5=\", \"79+5=\", \"38+7=\", \"43-16=\"),\n    @(\"77+4=\", \"42+19=\", \"93-7=\", \"62-35=\", \"88+3=\"),\n    @(\"75+17=\", \"45-16=\", \"83-37=\", \"17+68=\", \"16+16=\")\n)\n\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
